$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "245.43"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-0.60%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "27.19"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "2.91%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.107"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "0.69%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.05705"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "1.83%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "6.507"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "0.52%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.8192"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "0.71%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.8590"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "1.98%"
$ws.Range("B9").Value = "MandalaExchangeToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06935"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-0.67%"
$ws.Range("B10").Value = "BitrueCoin"
$ws.Range("C10").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.02844"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-0.70%"
$ws.Range("B11").Value = "BitMartToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.09395"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "0.08%"
$ws.Range("B12").Value = "BitForexToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.001529"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "0.60%"
$ws.Range("B13").Value = "CoinExToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.04045"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-13.04%"
$ws.Range("B14").Value = "One"
$ws.Range("C14").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0006017"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.62%"
$ws.Range("B15").Value = "TigerCash"
$ws.Range("C15").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.006212"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "0.31%"
$ws.Range("B16").Value = "LEO"
$ws.Range("C16").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.511"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-2.66%"
$ws.Range("B17").Value = "GateToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.008"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-0.26%"
$ws.Range("B18").Value = "BTSEToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.316"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "12.70%"
$ws.Range("B19").Value = "BitpandaEcosystemToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3165"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "1.24%"
$ws.Range("B20").Value = "WazirX"
$ws.Range("C20").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.1331"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-0.69%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.03225"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "1.02%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.1274"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-1.79%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.574"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-4.45%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.1373"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "1.74%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.001216"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-2.41%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.004473"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-2.48%"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.00009895"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "3.09%"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0001448"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "3.62%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03726"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "1.72%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.005907"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "72.64%"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-21.59%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002299"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-13.56%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.009220"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "11.45%"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-4.35%"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-0.04%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.1009"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.002505"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-3.81%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00002099"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-0.04%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0001999"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.04%"
